# Add 2022-Q1 data:
#  - the current "总计" (totals) sheet becomes "2022-Q1" and is filled with the
#    quarter's fund-holding detail rows
#  - a brand-new "总计" sheet is appended after it, carrying the same summary
#    table as before plus one new row for 2022-Q1

$wb = $excel.ActiveWorkbook

$q1Sheet = $wb.Worksheets.Item("总计")

# Grab a couple of already-styled cells from this sheet before we touch
# anything - we reuse their formatting (bold / border / centered) for the
# new header + index cells instead of re-building the style by hand.
$headerFmtSrc = $q1Sheet.Range("B1")
$indexFmtSrc = $q1Sheet.Range("A2")

# Rename the existing totals sheet to the new quarter, then append a fresh
# "总计" sheet right after it so the tab order ends up
# 2021-Q3, 2021-Q4, 2022-Q1, 总计
$q1Sheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Add($null, $q1Sheet)
$totalSheet.Name = "总计"

# ---------------------------------------------------------------------
# 2022-Q1 sheet: fund holding detail
# ---------------------------------------------------------------------

$q1Sheet.Cells.ClearContents()

$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $q1Headers.Length; $c++) {
    $q1Sheet.Cells.Item(1, $c + 2).Value = $q1Headers[$c]
}

# rows = idx, code, name, fundSize, stockPosition, positionPct, marketValue, rank, marketValueIsNumericZero
$q1Data = @(
    ,@(0, "003567", "华夏行业景气混合", "112.21", "91.63", "4.43", "4.9709", 2, $false)
    ,@(1, "519732", "交银定期支付双息平衡混合", "40.83", "67.67", "2.07", "0.8452", 5, $false)
    ,@(2, "519110", "浦银安盛价值成长混合A", "8.04", "88.87", "4.45", "0.3578", 2, $false)
    ,@(3, "000601", "华宝创新优选混合", "12.99", "87.56", "2.64", "0.3429", 10, $false)
    ,@(4, "519170", "浦银安盛增长动力灵活配置混合", "8.12", "85.61", "2.97", "0.2412", 5, $false)
    ,@(5, "001749", "招商中国机遇股票", "4.57", "94.84", "3.99", "0.1823", 6, $false)
    ,@(6, "007731", "民生加银持续成长混合A", "2.62", "93.83", "6.76", "0.1771", 2, $false)
    ,@(7, "519113", "浦银安盛精致生活混合", "2.09", "90.20", "5.25", "0.1097", 2, $false)
    ,@(8, "519120", "浦银安盛新兴产业混合", "2.21", "90.11", "4.77", "0.1054", 2, $false)
    ,@(9, "002103", "招商康泰灵活配置混合", "1.95", "39.68", "4.02", "0.0784", 3, $false)
    ,@(10, "350002", "天治低碳经济灵活配置混合", "0.76", "65.23", "5.74", "0.0436", 5, $false)
    ,@(11, "011599", "国联安匠心科技1个月滚动持有混合", "0.71", "91.41", "5.73", "0.0407", 7, $false)
    ,@(12, "007732", "民生加银持续成长混合C", "0.14", "93.83", "6.76", "0.0095", 2, $false)
    ,@(13, "960031", "浦银安盛价值成长混合H", "0.00", "88.87", "4.45", 0, 2, $true)
    ,@(14, "014011", "浦银安盛价值成长混合C", "0.00", "88.87", "4.45", 0, 2, $true)
)

# Text-ish columns must keep their number format as Text ("@") so values
# like fund codes ("003567") or percentages ("4.45") round-trip verbatim
# instead of being auto-coerced into numbers by Excel. The last two rows'
# market value rounds down to plain 0, which stays a real number.
$q1Sheet.Range("B2:B16").NumberFormat = "@"
$q1Sheet.Range("D2:F16").NumberFormat = "@"
$q1Sheet.Range("G2:G14").NumberFormat = "@"

foreach ($row in $q1Data) {
    $r = $row[0] + 2
    $q1Sheet.Cells.Item($r, 1).Value = $row[0]
    $q1Sheet.Cells.Item($r, 2).Value = $row[1]
    $q1Sheet.Cells.Item($r, 3).Value = $row[2]
    $q1Sheet.Cells.Item($r, 4).Value = $row[3]
    $q1Sheet.Cells.Item($r, 5).Value = $row[4]
    $q1Sheet.Cells.Item($r, 6).Value = $row[5]
    $q1Sheet.Cells.Item($r, 7).Value = $row[6]
    $q1Sheet.Cells.Item($r, 8).Value = $row[7]
}

# Re-apply the bold/border/center formatting that belongs on the header row
# and on the row-index column (column A).
$headerFmtSrc.Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)
$indexFmtSrc.Copy()
$q1Sheet.Range("A2:A16").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 总计 (totals) sheet: same summary table, with a new row for 2022-Q1
# ---------------------------------------------------------------------

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($c = 0; $c -lt $totalHeaders.Length; $c++) {
    $totalSheet.Cells.Item(1, $c + 2).Value = $totalHeaders[$c]
}

$totalData = @(
    ,@(0, "2022-Q1", 15, 7.5)
    ,@(1, "2021-Q4", 1, 4.61)
    ,@(2, "2021-Q3", 1, 2.08)
)

foreach ($row in $totalData) {
    $r = $row[0] + 2
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
}

$headerFmtSrc.Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)
$indexFmtSrc.Copy()
$totalSheet.Range("A2:A4").PasteSpecial(-4122)
